# Adds WAT106 and WAT107 test cases to the "Test Cases" sheet, matching
# the commit "Adds WAT107 and related changes".
#
# Two new rows are appended after the existing last row (103):
#   Row 104: WAT106 / WAT-330 / Verify that the Organization list provided
#            should be arranged alphabetically. / Y
#   Row 105: WAT107 / WAT-327 / Verify that the Country list provided
#            should be arranged alphabetically. / Y

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Seed the two new rows with the same look & feel (borders/wrap) as the
# existing table by copying the formatting of the last populated row (103)
# down onto the two freshly appended rows before filling in their values.
$ws.Range("A103:E103").Copy()
$ws.Range("A104:E105").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Values are entered column-by-column (A104, A105, B104, B105, C104, C105,
# ...) so new shared-string entries land in the same order the author's
# Excel session produced them in.
$ws.Cells.Item(104, 1).Value = "WAT106"
$ws.Cells.Item(105, 1).Value = "WAT107"

$ws.Cells.Item(104, 2).Value = "WAT-330"
$ws.Cells.Item(105, 2).Value = "WAT-327"

$ws.Cells.Item(104, 3).Value = "Verify that the Organization list provided should be arranged alphabetically."
$ws.Cells.Item(105, 3).Value = "Verify that the Country list provided should be arranged alphabetically."

$ws.Cells.Item(104, 4).Value = "Y"
$ws.Cells.Item(105, 4).Value = "Y"

# Leave the selection where the author last left it.
$ws.Range("C107").Select()
